$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the original sheet and add the three new sheets in order
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Prediction Info"

# ------------------------------------------------------------------
# 2. Capture the original "Sales vs PO" data (A=ds, B=y, C=PO_Requested_Qty)
#    before it gets overwritten.
# ------------------------------------------------------------------
$oldDs = @()
$oldY = @()
$oldPo = @()
for ($r = 2; $r -le 26; $r++) {
    $oldDs += $ws1.Cells.Item($r, 1).Value2
    $oldY += $ws1.Cells.Item($r, 2).Value2
    $oldPo += $ws1.Cells.Item($r, 3).Value2
}

# ------------------------------------------------------------------
# 3. Build the new "Sales vs PO" layout:
#    A=ds (shifted +6 days), B=y, C=Order Week (old ds), D=PO_Requested_Qty (old PO)
# ------------------------------------------------------------------

# Give the new column D the same bold header style as the other headers
$ws1.Cells.Item(1, 1).Copy($ws1.Cells.Item(1, 4))
# Give the new column C the same bold header style as the other headers
$ws1.Cells.Item(1, 1).Copy($ws1.Cells.Item(1, 3))

$ws1.Cells.Item(1, 1).Value2 = "ds"
$ws1.Cells.Item(1, 2).Value2 = "y"
$ws1.Cells.Item(1, 3).Value2 = "Order Week"
$ws1.Cells.Item(1, 4).Value2 = "PO_Requested_Qty"

# Give column C the same date number-format as column A
for ($r = 2; $r -le 26; $r++) {
    $ws1.Cells.Item($r, 1).Copy($ws1.Cells.Item($r, 3))
}

for ($r = 2; $r -le 26; $r++) {
    $i = $r - 2
    $ws1.Cells.Item($r, 1).Value2 = $oldDs[$i] + 6
    $ws1.Cells.Item($r, 2).Value2 = $oldY[$i]
    $ws1.Cells.Item($r, 3).Value2 = $oldDs[$i]
    $ws1.Cells.Item($r, 4).Value2 = $oldPo[$i]
}

# ------------------------------------------------------------------
# 4. "Weekly Growth" sheet: ds | PO_Requested_Qty | Growth%
# ------------------------------------------------------------------
$ws1.Cells.Item(1, 1).Copy($ws2.Cells.Item(1, 1))
$ws1.Cells.Item(1, 1).Copy($ws2.Cells.Item(1, 2))
$ws1.Cells.Item(1, 1).Copy($ws2.Cells.Item(1, 3))

$ws2.Cells.Item(1, 1).Value2 = "ds"
$ws2.Cells.Item(1, 2).Value2 = "PO_Requested_Qty"
$ws2.Cells.Item(1, 3).Value2 = "Growth%"

$growthDs = @(45488, 45551, 45579, 45586, 45593)
$growthPo = @(80, 256, 16, 336, 192)
$growthPct = @(0, 220, -93.75, 2000, -42.85714285714286)

for ($r = 2; $r -le 6; $r++) {
    $ws1.Cells.Item(2, 1).Copy($ws2.Cells.Item($r, 1))
}

for ($r = 2; $r -le 6; $r++) {
    $i = $r - 2
    $ws2.Cells.Item($r, 1).Value2 = $growthDs[$i]
    $ws2.Cells.Item($r, 2).Value2 = $growthPo[$i]
    $ws2.Cells.Item($r, 3).Value2 = $growthPct[$i]
}

# ------------------------------------------------------------------
# 5. "Volume Insights" sheet: Total / Average / Max / Min PO quantity
# ------------------------------------------------------------------
$ws1.Cells.Item(1, 1).Copy($ws3.Cells.Item(1, 1))
$ws1.Cells.Item(1, 1).Copy($ws3.Cells.Item(1, 2))
$ws1.Cells.Item(1, 1).Copy($ws3.Cells.Item(1, 3))
$ws1.Cells.Item(1, 1).Copy($ws3.Cells.Item(1, 4))

$ws3.Cells.Item(1, 1).Value2 = "Total_PO_Quantity"
$ws3.Cells.Item(1, 2).Value2 = "Average_PO_Quantity"
$ws3.Cells.Item(1, 3).Value2 = "Max_PO_Quantity"
$ws3.Cells.Item(1, 4).Value2 = "Min_PO_Quantity"

$ws3.Cells.Item(2, 1).Value2 = 880
$ws3.Cells.Item(2, 2).Value2 = 176
$ws3.Cells.Item(2, 3).Value2 = 336
$ws3.Cells.Item(2, 4).Value2 = 16

# ------------------------------------------------------------------
# 6. "Prediction Info" sheet: Predicted_Next_Week_PO_Quantity
# ------------------------------------------------------------------
$ws1.Cells.Item(1, 1).Copy($ws4.Cells.Item(1, 1))
$ws4.Cells.Item(1, 1).Value2 = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2, 1).Value2 = 267.2
